$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New label above the existing (first) measurement table ---
$ws.Range("D3").Value = "Prototyp"

# --- Second table title ---
$ws.Range("D11").Value = "Finální systém"

# --- Second table header row (same headers as the first table) ---
$ws.Range("D12").Value = "PCAP size (Bytes)"
$ws.Range("E12").Value = "Count"
$ws.Range("F12").Value = "AVG packet size"
$ws.Range("G12").Value = "Time (s)"
$ws.Range("H12").Value = "Packet/s"
$ws.Range("J12").Value = "AVG time (s)"
$ws.Range("K12").Value = "AVG packet size (Bytes)"

# --- Second table data rows 13-17 (same PCAP/count data as rows 5-9, different interval G) ---
$ws.Range("D13").Value = 91340
$ws.Range("E13").Value = 382
$ws.Range("G13").Value = 0.5

$ws.Range("D14").Value = 1955172
$ws.Range("E14").Value = 1360
$ws.Range("G14").Value = 0.5

$ws.Range("D15").Value = 412254
$ws.Range("E15").Value = 1879
$ws.Range("G15").Value = 0.5

$ws.Range("D16").Value = 420869
$ws.Range("E16").Value = 2263
$ws.Range("G16").Value = 1

$ws.Range("D17").Value = 449234
$ws.Range("E17").Value = 2563
$ws.Range("G17").Value = 1.5

# F column: AVG packet size = D/E, row 13 gets its own formula, 14:17 entered as one
# block so the engine stores them as a shared formula (mirrors existing F6:F9 block)
$ws.Range("F13").Formula = "=D13/E13"
$ws.Range("F14:F17").Formula = "=D14/E14"

# H column: Packet/s = E/G (one formula per row, matches the first table's layout)
$ws.Range("H13").Formula = "=E13/G13"
$ws.Range("H14").Formula = "=E14/G14"
$ws.Range("H15").Formula = "=E15/G15"
$ws.Range("H16").Formula = "=E16/G16"
$ws.Range("H17").Formula = "=E17/G17"

# J13/K13: averages over the new block
$ws.Range("J13").Formula = "=AVERAGEA(H13:H17)"
$ws.Range("K13").Formula = "=AVERAGEA(F13:F17)"

# number format for the byte-size column (match the first table's "D" column style) -
# applied last so it only lands on the D cells, not any formula cells created above
$ws.Range("D13:D17").NumberFormat = "#,##0"

# Move the active selection like the saved workbook shows
$ws.Range("K21").Select() | Out-Null

Write-Host "edit applied"
